$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 47
$ws1.Range("F5").Value = 735
$ws1.Range("F6").Value = 2275
$ws1.Range("F8").Value = 1709
$ws1.Range("F9").Value = 2882
$ws1.Range("F10").Value = 159
$ws1.Range("F11").Value = 4285
$ws1.Range("F12").Value = 370
$ws1.Range("F17").Value = 4
$ws1.Range("F21").Value = 293
$ws1.Range("F22").Value = 4135
$ws1.Range("F24").Value = 3619
$ws1.Range("F25").Value = 1124
$ws1.Range("F26").Value = 209
$ws1.Range("F27").Value = 533
$ws1.Range("F28").Value = 4364
$ws1.Range("F29").Value = 82
$ws1.Range("F30").Value = 476
$ws1.Range("F31").Value = 509
$ws1.Range("F32").Value = 449

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 5

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 47
$ws4.Range("F7").Value = 735
$ws4.Range("F8").Value = 2275
$ws4.Range("F10").Value = 1709
$ws4.Range("F12").Value = 2882
$ws4.Range("F13").Value = 159
$ws4.Range("F14").Value = 4285
$ws4.Range("F15").Value = 370
$ws4.Range("F20").Value = 4
$ws4.Range("F25").Value = 293
$ws4.Range("F26").Value = 4135
$ws4.Range("F28").Value = 3619
$ws4.Range("F29").Value = 1124
$ws4.Range("F30").Value = 209
$ws4.Range("F31").Value = 533
$ws4.Range("F32").Value = 4364
$ws4.Range("F33").Value = 82
$ws4.Range("F34").Value = 476
$ws4.Range("F35").Value = 509
$ws4.Range("F36").Value = 449
$ws4.Range("F37").Value = 5
